$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.676.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.422.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.21"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.84"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.516"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.82%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +10.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.50"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.84"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.801.42"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.406.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.511.95"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.71"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.19"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.55"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.56"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.53"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +18.05%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +11.34%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0777"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.97%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.50"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.13"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.09"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.950.56"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +8.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.70"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +12.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.31"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.43%  "
